# Update "想去人数" (F column) figures across the 展览/演出/全部类型 sheets
# to reflect the regenerated gh-pages output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# ---- 展览 sheet ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 794
$ws.Range("F3").Value = 531
$ws.Range("F4").Value = 280
$ws.Range("F5").Value = 488
$ws.Range("F6").Value = 1131
$ws.Range("F8").Value = 31
$ws.Range("F9").Value = 112
$ws.Range("F11").Value = 1142
$ws.Range("F14").Value = 782
$ws.Range("F15").Value = 812
$ws.Range("F16").Value = 183
$ws.Range("F17").Value = 44
$ws.Range("F18").Value = 65
$ws.Range("F19").Value = 665
$ws.Range("F20").Value = 185
$ws.Range("F21").Value = 1713
$ws.Range("F22").Value = 2293
$ws.Range("F23").Value = 628
$ws.Range("F24").Value = 65
$ws.Range("F25").Value = 1884
$ws.Range("F26").Value = 319
$ws.Range("F27").Value = 2745
$ws.Range("F28").Value = 504
$ws.Range("F30").Value = 679
$ws.Range("F31").Value = 130
$ws.Range("F32").Value = 97
$ws.Range("F34").Value = 954
$ws.Range("F35").Value = 1681
$ws.Range("F36").Value = 322
$ws.Range("F39").Value = 153
$ws.Range("F41").Value = 153
$ws.Range("F42").Value = 9

# ---- 演出 sheet ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 135
$ws.Range("F4").Value = 10
$ws.Range("F9").Value = 8
$ws.Range("F12").Value = 70

# ---- 全部类型 sheet ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 794
$ws.Range("F4").Value = 531
$ws.Range("F5").Value = 280
$ws.Range("F6").Value = 488
$ws.Range("F7").Value = 1131
$ws.Range("F9").Value = 31
$ws.Range("F10").Value = 112
$ws.Range("F12").Value = 1142
$ws.Range("F14").Value = 782
$ws.Range("F15").Value = 812
$ws.Range("F16").Value = 183
$ws.Range("F17").Value = 135
$ws.Range("F18").Value = 135
$ws.Range("F19").Value = 10
$ws.Range("F20").Value = 44
$ws.Range("F22").Value = 65
$ws.Range("F23").Value = 185
$ws.Range("F24").Value = 1713
$ws.Range("F25").Value = 2293
$ws.Range("F26").Value = 628
$ws.Range("F27").Value = 65
$ws.Range("F30").Value = 2745
$ws.Range("F31").Value = 504
$ws.Range("F33").Value = 8
$ws.Range("F37").Value = 70
$ws.Range("F38").Value = 679
$ws.Range("F39").Value = 130
$ws.Range("F40").Value = 97
$ws.Range("F42").Value = 954
$ws.Range("F43").Value = 1681
$ws.Range("F45").Value = 322
$ws.Range("F47").Value = 153
$ws.Range("F49").Value = 153
